{"js": "// Update the two-digit-divided-by-one-digit division problems in the\n// worksheet table. Each \"old=>new\" pair is unique within the document, so\n// a simple search-and-replace per pair is sufficient and unambiguous.\nconst replacements = [\n  [\"36\u00f78=\", \"35\u00f75=\"],\n  [\"67\u00f72=\", \"16\u00f76=\"],\n  [\"11\u00f74=\", \"20\u00f74=\"],\n  [\"75\u00f76=\", \"44\u00f72=\"],\n  [\"26\u00f79=\", \"13\u00f72=\"],\n  [\"18\u00f72=\", \"61\u00f73=\"],\n  [\"86\u00f77=\", \"66\u00f79=\"],\n  [\"21\u00f72=\", \"79\u00f79=\"],\n  [\"43\u00f74=\", \"38\u00f72=\"],\n  [\"73\u00f72=\", \"17\u00f76=\"],\n  [\"98\u00f78=\", \"90\u00f79=\"],\n  [\"55\u00f72=\", \"47\u00f78=\"],\n  [\"29\u00f79=\", \"94\u00f74=\"],\n  [\"33\u00f78=\", \"60\u00f76=\"],\n  [\"21\u00f74=\", \"61\u00f79=\"],\n  [\"93\u00f79=\", \"66\u00f78=\"],\n  [\"87\u00f79=\", \"40\u00f77=\"],\n  [\"85\u00f78=\", \"64\u00f78=\"],\n  [\"60\u00f73=\", \"49\u00f75=\"],\n  [\"17\u00f72=\", \"31\u00f74=\"],\n  [\"83\u00f73=\", \"33\u00f77=\"],\n  [\"41\u00f73=\", \"76\u00f72=\"],\n  [\"24\u00f77=\", \"50\u00f74=\"],\n  [\"87\u00f77=\", \"15\u00f77=\"],\n  [\"23\u00f79=\", \"61\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-divided-by-one-digit division problems in the\n# worksheet table. Each \"old=>new\" pair is unique within the document, so\n# a Find/Replace-All per pair is sufficient and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"36\u00f78=\", \"35\u00f75=\"),\n    @(\"67\u00f72=\", \"16\u00f76=\"),\n    @(\"11\u00f74=\", \"20\u00f74=\"),\n    @(\"75\u00f76=\", \"44\u00f72=\"),\n    @(\"26\u00f79=\", \"13\u00f72=\"),\n    @(\"18\u00f72=\", \"61\u00f73=\"),\n    @(\"86\u00f77=\", \"66\u00f79=\"),\n    @(\"21\u00f72=\", \"79\u00f79=\"),\n    @(\"43\u00f74=\", \"38\u00f72=\"),\n    @(\"73\u00f72=\", \"17\u00f76=\"),\n    @(\"98\u00f78=\", \"90\u00f79=\"),\n    @(\"55\u00f72=\", \"47\u00f78=\"),\n    @(\"29\u00f79=\", \"94\u00f74=\"),\n    @(\"33\u00f78=\", \"60\u00f76=\"),\n    @(\"21\u00f74=\", \"61\u00f79=\"),\n    @(\"93\u00f79=\", \"66\u00f78=\"),\n    @(\"87\u00f79=\", \"40\u00f77=\"),\n    @(\"85\u00f78=\", \"64\u00f78=\"),\n    @(\"60\u00f73=\", \"49\u00f75=\"),\n    @(\"17\u00f72=\", \"31\u00f74=\"),\n    @(\"83\u00f73=\", \"33\u00f77=\"),\n    @(\"41\u00f73=\", \"76\u00f72=\"),\n    @(\"24\u00f77=\", \"50\u00f74=\"),\n    @(\"87\u00f77=\", \"15\u00f77=\"),\n    @(\"23\u00f79=\", \"61\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    # 2 = wdReplaceAll\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
